$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 190, pushing existing rows 190..207 down to 191..208
$ws.Rows.Item(190).Insert()

# Populate the newly inserted row 190 with the new weekly data point
$ws.Cells.Item(190, 1).Value = 8
$ws.Cells.Item(190, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(190, 3).Value = "Coquimbo"
$ws.Cells.Item(190, 4).Value = 45142
$ws.Cells.Item(190, 5).Value = 4
$ws.Cells.Item(190, 6).Value = 100112052
$ws.Cells.Item(190, 7).Value = "Albahaca"
$ws.Cells.Item(190, 8).Value = "Sin especificar"
$ws.Cells.Item(190, 9).Value = "Primera"
$ws.Cells.Item(190, 10).Value = 800
$ws.Cells.Item(190, 11).Value = 2800
$ws.Cells.Item(190, 12).Value = 3000
$ws.Cells.Item(190, 13).Value = 2900
$ws.Cells.Item(190, 14).Value = "`$/paquete"
$ws.Cells.Item(190, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(190, 16).Value = 2900
$ws.Cells.Item(190, 17).Value = 1
$ws.Cells.Item(190, 18).Value = "Hortaliza"
